$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: values are prefixed with a leading apostrophe so Excel stores them
# as text (matching the original inlineStr/text cells) instead of coercing
# numeric-looking strings (e.g. "97.10") into floating point numbers that
# would lose trailing zeros / the original textual representation.

$ws.Range("D2").Value = "'35.467.64"
$ws.Range("E2").Value = "'  +0.77%  "
$ws.Range("D3").Value = "'1.923.33"
$ws.Range("E3").Value = "'  +1.80%  "
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("D5").Value = "'0.729"
$ws.Range("E5").Value = "'  +11.65%  "
$ws.Range("D6").Value = "'254.41"
$ws.Range("E6").Value = "'  +4.89%  "
$ws.Range("E7").Value = "'  +0.04%  "
$ws.Range("D8").Value = "'40.94"
$ws.Range("E8").Value = "'  -0.34%  "
$ws.Range("E9").Value = "'  +2.93%  "
$ws.Range("E10").Value = "'  +5.06%  "
$ws.Range("D11").Value = "'0.0749"
$ws.Range("E11").Value = "'  +6.13%  "
$ws.Range("D12").Value = "'0.0998"
$ws.Range("E12").Value = "'  +0.21%  "
$ws.Range("D13").Value = "'2.203.76"
$ws.Range("E13").Value = "'  +1.72%  "
$ws.Range("D14").Value = "'12.79"
$ws.Range("E14").Value = "'  +7.82%  "
$ws.Range("E15").Value = "'  +4.52%  "
$ws.Range("D16").Value = "'1.926.79"
$ws.Range("E16").Value = "'  +1.99%  "
$ws.Range("E17").Value = "'  +1.99%  "
$ws.Range("D18").Value = "'35.470.30"
$ws.Range("D19").Value = "'74.53"
$ws.Range("E19").Value = "'  +4.93%  "
$ws.Range("D20").Value = "'0.0₃0839"
$ws.Range("E20").Value = "'  +3.70%  "
$ws.Range("D21").Value = "'243.69"
$ws.Range("E21").Value = "'  +1.49%  "
$ws.Range("D22").Value = "'13.04"
$ws.Range("E22").Value = "'  +5.63%  "
$ws.Range("D23").Value = "'5.12"
$ws.Range("E23").Value = "'  +8.63%  "
$ws.Range("E24").Value = "'  +0.01%  "
$ws.Range("E25").Value = "'  +2.57%  "
$ws.Range("D26").Value = "'2.41"
$ws.Range("E26").Value = "'  -0.82%  "
$ws.Range("D27").Value = "'167.66"
$ws.Range("E27").Value = "'  -1.31%  "
$ws.Range("D28").Value = "'8.65"
$ws.Range("E28").Value = "'  +3.11%  "
$ws.Range("E29").Value = "'  +6.72%  "
$ws.Range("D30").Value = "'18.77"
$ws.Range("E30").Value = "'  +3.39%  "
$ws.Range("D31").Value = "'4.128.27"
$ws.Range("E31").Value = "'  +19.44%  "
$ws.Range("E32").Value = "'  +7.48%  "
$ws.Range("B33").Value = "'WEMIXToken"
$ws.Range("C33").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").Value = "'1.99"
$ws.Range("E33").Value = "'  +14.91%  "
$ws.Range("B34").Value = "'TrustWalletToken"
$ws.Range("C34").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D34").Value = "'1.64"
$ws.Range("E34").Value = "'  +24.99%  "
$ws.Range("D35").Value = "'0.0583"
$ws.Range("E35").Value = "'  +4.45%  "
$ws.Range("E36").Value = "'  +4.02%  "
$ws.Range("E37").Value = "'  -0.03%  "
$ws.Range("E38").Value = "'  -1.78%  "
$ws.Range("E39").Value = "'  +1.25%  "
$ws.Range("D40").Value = "'17.49"
$ws.Range("E40").Value = "'  +10.15%  "
$ws.Range("D41").Value = "'97.10"
$ws.Range("E41").Value = "'  +9.48%  "
$ws.Range("E42").Value = "'  +3.97%  "
$ws.Range("B43").Value = "'VeChain"
$ws.Range("C43").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0211"
$ws.Range("E43").Value = "'  +1.73%  "
$ws.Range("B44").Value = "'Kaspa"
$ws.Range("C44").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").Value = "'0.0657"
$ws.Range("E44").Value = "'  +3.01%  "
$ws.Range("D45").Value = "'1.347.61"
$ws.Range("E45").Value = "'  +0.99%  "
$ws.Range("E46").Value = "'  +4.83%  "
$ws.Range("B47").Value = "'HuobiToken"
$ws.Range("C47").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D47").Value = "'2.43"
$ws.Range("E47").Value = "'  +1.21%  "
$ws.Range("B48").Value = "'FraxShare"
$ws.Range("C48").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'6.80"
$ws.Range("E48").Value = "'  +4.43%  "
$ws.Range("E49").Value = "'  +0.82%  "
$ws.Range("D50").Value = "'45.27"
$ws.Range("E50").Value = "'  -5.39%  "
$ws.Range("D51").Value = "'11.95"
$ws.Range("E51").Value = "'  +7.06%  "
